# batch_create_jdbc_datasource-en.xlsx : fix(cdk): batch create template
#
# Changes applied:
#  1. Column header in F2 renamed "SecretARN" -> "SecretID"
#  2. Matching bold label inside the A1 instructions cell renamed
#     "SecretARN" -> "SecretID" (only the standalone bold header mention -
#     the two "(... ignored if SecretARN is filled)" remarks stay as-is)
#  3. Selection moved from A1:K1 to the single cell D4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header cell text -------------------------------------------------
$ws.Range("F2").Value = "SecretID"

# --- 2. Rich-text instructions cell --------------------------------------
# Capture each existing run's text + formatting BEFORE touching the text,
# since re-assigning .Characters(...).Text collapses every run into one
# plain run. We rebuild the exact same run boundaries afterwards (offsets
# after the target run shift by -1 because "SecretID" is one character
# shorter than "SecretARN"). Plain arrays are used throughout (rather than
# a collection of custom objects) since that indexes/reads back reliably.

$cell = $ws.Range("A1")
$full = $cell.Value2

$starts = @()
$lens = @()
$bolds = @()
$sizes = @()
$names = @()
$texts = @()

$pos = 0
while ($pos -lt $full.Length) {
    $run = $cell.Characters($pos + 1, 1)
    $bold = $run.Font.Bold
    $size = $run.Font.Size
    $name = $run.Font.Name
    $len = 1
    while ($pos + $len -lt $full.Length) {
        $next = $cell.Characters($pos + $len + 1, 1)
        if ($next.Font.Bold -ne $bold -or $next.Font.Size -ne $size -or $next.Font.Name -ne $name) {
            break
        }
        $len = $len + 1
    }
    $starts += $pos
    $lens += $len
    $bolds += $bold
    $sizes += $size
    $names += $name
    $texts += $full.Substring($pos, $len)
    $pos = $pos + $len
}

# Find the standalone bold "SecretARN" run (the column-style label, not the
# inline mentions inside the Username/Password remarks) and rename it.
$targetIndex = -1
for ($i = 0; $i -lt $starts.Length; $i++) {
    if ($bolds[$i] -and $texts[$i] -eq "SecretARN") {
        $targetIndex = $i
    }
}

if ($targetIndex -ge 0) {
    $tStart = $starts[$targetIndex]
    $tLen = $lens[$targetIndex]
    $newFull = $full.Substring(0, $tStart) + "SecretID" + $full.Substring($tStart + $tLen)

    # Rewriting .Text collapses formatting, so do it once up-front...
    $cell.Characters(1, $full.Length).Text = $newFull

    # ...then walk the recorded runs again (using the renamed run's new,
    # shorter length) and re-apply each run's captured font formatting.
    $cursor = 0
    for ($i = 0; $i -lt $starts.Length; $i++) {
        $len = $lens[$i]
        if ($i -eq $targetIndex) {
            $len = 8   # "SecretID".Length
        }
        $rng = $cell.Characters($cursor + 1, $len)
        $rng.Font.Bold = $bolds[$i]
        $rng.Font.Size = $sizes[$i]
        $rng.Font.Name = $names[$i]
        $cursor = $cursor + $len
    }
}

# --- 3. Selection ----------------------------------------------------------
$ws.Range("D4").Select()
